# TC02_Canine_Filter_Diagnosis-LipLingual.xlsx
# "updated statbar xpaths & diagnosis testcases"
#
# The "startup" sheet gains a new second column ("StatQuery") holding a
# companion Neo4j query (a stats/count query) right next to the existing
# "query" column. The columns that used to be B (dbExcel) and C (WebExcel)
# shift right to become C and D.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at B; this shifts the old B (dbExcel) and
# old C (WebExcel) columns one position to the right (-> C, D).
[void]$ws.Columns.Item(2).Insert()

# Header for the newly inserted column.
$ws.Range("B1").Value = "StatQuery"

# The new stats/count companion query for the diagnosis filter query in A2.
$ws.Range("B2").Value = "MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE diag.disease_term IN['Lip and oropharyngeal neoplasms malignant :: Melanoma-lingual']  OPTIONAL MATCH (f:file)-[*]->(c), (samp:sample)-[*]->(c) WITH DISTINCT c AS c, p, s, demo, diag, f, samp RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(samp)) as number_of_sample , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(s.clinical_study_designation)) as number_of_study"

# Match the wrap-text formatting already used by the neighboring A2 query cell.
$ws.Range("B2").WrapText = $true

# The new column should be exactly as wide as column A (the other query column).
$ws.Columns.Item(2).ColumnWidth = $ws.Columns.Item(1).ColumnWidth

# Reset the view: select B2 only (instead of the old B2:B6) and scroll back
# to show row 1 (instead of being scrolled to row 2).
[void]$ws.Range("A1").Select()
[void]$ws.Range("B2").Select()
